$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (B2, C2 string change; G2:T2 numeric updates)
# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Epcam"
$ws.Range("C2").Value2 = "Epcam"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.01553533333333333
$ws.Range("H2").Value2 = 0.046606
$ws.Range("I2").Value2 = 0.1067732734624681
$ws.Range("J2").Value2 = 0.1067732734624681
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 0.01553533333333333
$ws.Range("N2").Value2 = 0.046606
$ws.Range("O2").Value2 = 0.1067732734624681
$ws.Range("P2").Value2 = 0.1067732734624681
$ws.Range("Q2").Value2 = 0.0002413465817777778
$ws.Range("R2").Value2 = 0.002172119236
$ws.Range("S2").Value2 = 0.01140053192589099
$ws.Range("T2").Value2 = 0.01140053192589099

# Row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Epcam"
$ws.Range("C3").Value2 = "Epcam"
$ws.Range("D3").Value2 = "Inflammatory-Mac"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.01553533333333333
$ws.Range("H3").Value2 = 0.046606
$ws.Range("I3").Value2 = 0.1067732734624681
$ws.Range("J3").Value2 = 0.1067732734624681
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0.3333333333333333
$ws.Range("M3").Value2 = 0.1140266666666667
$ws.Range("N3").Value2 = 0.34208
$ws.Range("O3").Value2 = 0.783697407759539
$ws.Range("P3").Value2 = 0.7836974077595391
$ws.Range("Q3").Value2 = 0.001771442275555556
$ws.Range("R3").Value2 = 0.01594298048
$ws.Range("S3").Value2 = 0.0836779376305366
$ws.Range("T3").Value2 = 0.08367793763053662

# Row 4
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Epcam"
$ws.Range("C4").Value2 = "Epcam"
$ws.Range("D4").Value2 = "Resolving-Mac"
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.01553533333333333
$ws.Range("H4").Value2 = 0.046606
$ws.Range("I4").Value2 = 0.1067732734624681
$ws.Range("J4").Value2 = 0.1067732734624681
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.01593633333333333
$ws.Range("N4").Value2 = 0.047809
$ws.Range("O4").Value2 = 0.1095293187779929
$ws.Range("P4").Value2 = 0.1095293187779929
$ws.Range("Q4").Value2 = 0.0002475762504444444
$ws.Range("R4").Value2 = 0.002228186254
$ws.Range("S4").Value2 = 0.01169480390604047
$ws.Range("T4").Value2 = 0.01169480390604047

# Row 5
$ws.Range("A5").Value2 = "Inflammatory-Mac"
$ws.Range("B5").Value2 = "Epcam"
$ws.Range("C5").Value2 = "Epcam"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.1140266666666667
$ws.Range("H5").Value2 = 0.34208
$ws.Range("I5").Value2 = 0.783697407759539
$ws.Range("J5").Value2 = 0.7836974077595391
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.01553533333333333
$ws.Range("N5").Value2 = 0.046606
$ws.Range("O5").Value2 = 0.1067732734624681
$ws.Range("P5").Value2 = 0.1067732734624681
$ws.Range("Q5").Value2 = 0.001771442275555556
$ws.Range("R5").Value2 = 0.01594298048
$ws.Range("S5").Value2 = 0.0836779376305366
$ws.Range("T5").Value2 = 0.08367793763053662

# Row 6
$ws.Range("A6").Value2 = "Inflammatory-Mac"
$ws.Range("B6").Value2 = "Epcam"
$ws.Range("C6").Value2 = "Epcam"
$ws.Range("D6").Value2 = "Inflammatory-Mac"
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.1140266666666667
$ws.Range("H6").Value2 = 0.34208
$ws.Range("I6").Value2 = 0.783697407759539
$ws.Range("J6").Value2 = 0.7836974077595391
$ws.Range("K6").Value2 = 1
$ws.Range("L6").Value2 = 0.3333333333333333
$ws.Range("M6").Value2 = 0.1140266666666667
$ws.Range("N6").Value2 = 0.34208
$ws.Range("O6").Value2 = 0.783697407759539
$ws.Range("P6").Value2 = 0.7836974077595391
$ws.Range("Q6").Value2 = 0.01300208071111111
$ws.Range("R6").Value2 = 0.1170187264
$ws.Range("S6").Value2 = 0.6141816269290211
$ws.Range("T6").Value2 = 0.6141816269290213

# Row 7
$ws.Range("A7").Value2 = "Inflammatory-Mac"
$ws.Range("B7").Value2 = "Epcam"
$ws.Range("C7").Value2 = "Epcam"
$ws.Range("D7").Value2 = "Resolving-Mac"
$ws.Range("E7").Value2 = 1
$ws.Range("F7").Value2 = 0.3333333333333333
$ws.Range("G7").Value2 = 0.1140266666666667
$ws.Range("H7").Value2 = 0.34208
$ws.Range("I7").Value2 = 0.783697407759539
$ws.Range("J7").Value2 = 0.7836974077595391
$ws.Range("K7").Value2 = 1
$ws.Range("L7").Value2 = 0.3333333333333333
$ws.Range("M7").Value2 = 0.01593633333333333
$ws.Range("N7").Value2 = 0.047809
$ws.Range("O7").Value2 = 0.1095293187779929
$ws.Range("P7").Value2 = 0.1095293187779929
$ws.Range("Q7").Value2 = 0.001817166968888889
$ws.Range("R7").Value2 = 0.01635450272
$ws.Range("S7").Value2 = 0.08583784319998121
$ws.Range("T7").Value2 = 0.08583784319998122

# Row 8
$ws.Range("A8").Value2 = "Resolving-Mac"
$ws.Range("B8").Value2 = "Epcam"
$ws.Range("C8").Value2 = "Epcam"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 1
$ws.Range("F8").Value2 = 0.3333333333333333
$ws.Range("G8").Value2 = 0.01593633333333333
$ws.Range("H8").Value2 = 0.047809
$ws.Range("I8").Value2 = 0.1095293187779929
$ws.Range("J8").Value2 = 0.1095293187779929
$ws.Range("K8").Value2 = 1
$ws.Range("L8").Value2 = 0.3333333333333333
$ws.Range("M8").Value2 = 0.01553533333333333
$ws.Range("N8").Value2 = 0.046606
$ws.Range("O8").Value2 = 0.1067732734624681
$ws.Range("P8").Value2 = 0.1067732734624681
$ws.Range("Q8").Value2 = 0.0002475762504444444
$ws.Range("R8").Value2 = 0.002228186254
$ws.Range("S8").Value2 = 0.01169480390604047
$ws.Range("T8").Value2 = 0.01169480390604047

# Row 9
$ws.Range("A9").Value2 = "Resolving-Mac"
$ws.Range("B9").Value2 = "Epcam"
$ws.Range("C9").Value2 = "Epcam"
$ws.Range("D9").Value2 = "Inflammatory-Mac"
$ws.Range("E9").Value2 = 1
$ws.Range("F9").Value2 = 0.3333333333333333
$ws.Range("G9").Value2 = 0.01593633333333333
$ws.Range("H9").Value2 = 0.047809
$ws.Range("I9").Value2 = 0.1095293187779929
$ws.Range("J9").Value2 = 0.1095293187779929
$ws.Range("K9").Value2 = 1
$ws.Range("L9").Value2 = 0.3333333333333333
$ws.Range("M9").Value2 = 0.1140266666666667
$ws.Range("N9").Value2 = 0.34208
$ws.Range("O9").Value2 = 0.783697407759539
$ws.Range("P9").Value2 = 0.7836974077595391
$ws.Range("Q9").Value2 = 0.001817166968888889
$ws.Range("R9").Value2 = 0.01635450272
$ws.Range("S9").Value2 = 0.08583784319998121
$ws.Range("T9").Value2 = 0.08583784319998122

# Row 10
$ws.Range("A10").Value2 = "Resolving-Mac"
$ws.Range("B10").Value2 = "Epcam"
$ws.Range("C10").Value2 = "Epcam"
$ws.Range("D10").Value2 = "Resolving-Mac"
$ws.Range("E10").Value2 = 1
$ws.Range("F10").Value2 = 0.3333333333333333
$ws.Range("G10").Value2 = 0.01593633333333333
$ws.Range("H10").Value2 = 0.047809
$ws.Range("I10").Value2 = 0.1095293187779929
$ws.Range("J10").Value2 = 0.1095293187779929
$ws.Range("K10").Value2 = 1
$ws.Range("L10").Value2 = 0.3333333333333333
$ws.Range("M10").Value2 = 0.01593633333333333
$ws.Range("N10").Value2 = 0.047809
$ws.Range("O10").Value2 = 0.1095293187779929
$ws.Range("P10").Value2 = 0.1095293187779929
$ws.Range("Q10").Value2 = 0.0002539667201111111
$ws.Range("R10").Value2 = 0.002285700481
$ws.Range("S10").Value2 = 0.01199667167197118
$ws.Range("T10").Value2 = 0.01199667167197118
